# "Generate Report for Handback"
#
# A new handback round completed for the file
# "4b4c079f-0cbc-42bd-a447-1df24f3b4675" in both the zh-cn and de-de
# target languages. The handoff/handback generation timestamps for that
# file's row need to be refreshed on each per-language sheet, and the
# roll-up "Latest HO Xliff Generate Date" on the Overview sheet needs to
# reflect the newest of those timestamps.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row for 4b4c079f...zh-cn.xlf ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H2").Value = "2016-08-21 12:53:12"   # Correspond Handoff Datetime
$wsZh.Range("K2").Value = "2016-08-21 12:53:29"   # Correspond Handback DateTime

# --- de-de sheet: row for 4b4c079f...de-de.xlf ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H2").Value = "2016-08-21 12:53:16"   # Correspond Handoff Datetime
$wsDe.Range("K2").Value = "2016-08-21 12:53:35"   # Correspond Handback DateTime

# --- Overview sheet: roll-up "Latest HO Xliff Generate Date" for the same file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-21 12:53:16"
